# Add a new "2020" column (N) to the indicator table, mirroring the
# formatting of the existing 2019 column (M) for each row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# xlPasteFormats = -4122 : paste only number format / font / borders / alignment
$xlPasteFormats = -4122

# --- Row 3 (thin bottom-border row under the header) : empty cell, just formatting
$ws.Range("M3").Copy()
$ws.Range("N3").PasteSpecial($xlPasteFormats)

# --- Row 4 (year header row) : 2020
$ws.Range("N4").Value = 2020
$ws.Range("M4").Copy()
$ws.Range("N4").PasteSpecial($xlPasteFormats)

# --- Row 5 (Kyrgyz Republic total)
$ws.Range("N5").Value = 95.9
$ws.Range("M5").Copy()
$ws.Range("N5").PasteSpecial($xlPasteFormats)

# --- Row 6 (Batken oblast)
$ws.Range("N6").Value = 96.2
$ws.Range("M6").Copy()
$ws.Range("N6").PasteSpecial($xlPasteFormats)
$ws.Range("N6").NumberFormat = "0.0"

# --- Row 7 (Djalal-Abad oblast)
$ws.Range("N7").Value = 62.3
$ws.Range("M7").Copy()
$ws.Range("N7").PasteSpecial($xlPasteFormats)
$ws.Range("N7").NumberFormat = "0.0"

# --- Row 8 (Ysyk-Kul oblast)
$ws.Range("N8").Value = 100
$ws.Range("M8").Copy()
$ws.Range("N8").PasteSpecial($xlPasteFormats)
$ws.Range("N8").NumberFormat = "0.0"

# --- Row 9 (Naryn oblast)
$ws.Range("N9").Value = 100
$ws.Range("M9").Copy()
$ws.Range("N9").PasteSpecial($xlPasteFormats)
$ws.Range("N9").NumberFormat = "0.0"

# --- Row 10 (Osh oblast) : "-"
$ws.Range("N10").Value = "-"
$ws.Range("M10").Copy()
$ws.Range("N10").PasteSpecial($xlPasteFormats)
$ws.Range("N10").NumberFormat = "0.0"

# --- Row 11 (Talas oblast)
$ws.Range("N11").Value = 100
$ws.Range("M11").Copy()
$ws.Range("N11").PasteSpecial($xlPasteFormats)
$ws.Range("N11").NumberFormat = "0.0"

# --- Row 12 (Chui oblast)
$ws.Range("N12").Value = 62.7
$ws.Range("M12").Copy()
$ws.Range("N12").PasteSpecial($xlPasteFormats)
$ws.Range("N12").NumberFormat = "0.0"

# --- Row 13 (Bishkek city)
$ws.Range("N13").Value = 100
$ws.Range("M13").Copy()
$ws.Range("N13").PasteSpecial($xlPasteFormats)
$ws.Range("N13").NumberFormat = "0.0"

# --- Row 14 (Osh city) : "-"
$ws.Range("N14").Value = "-"
$ws.Range("M14").Copy()
$ws.Range("N14").PasteSpecial($xlPasteFormats)
$ws.Range("N14").NumberFormat = "0.0"

# Clear the clipboard marquee and leave the new cell selected, matching
# the saved workbook's UI state.
$excel.CutCopyMode = $false
$ws.Range("N3").Select()
